$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2-34: update date serial from 45212 to 45221
$ws.Range("C2:C34").Value = 45221
